# conceptversie 0807 naar Teams
#
# Applies the edits captured in the commit diff:
#  - fixes a mojibake string ("ExponentiÃ«le" -> "Exponentiële") used on the
#    "A 2021" and "A 2020" sheets (H9)
#  - resets a batch of "herkans-baar?" (column O) dropdown cells that had
#    been filled in with "nee" back to their placeholder/true value on the
#    "A 2021", "A 2020", "A 2019" and "A 2018" sheets
#  - a couple of small numeric edits (I / N columns)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "A 2021"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("A 2021")

$ws.Range("H9").Value = "Hoofdstuk 4 (Machtsfuncties) + Hoofdstuk 5 (Exponentiële functies) + Vaardigheden"

$ws.Range("I6").Value = 1
$ws.Range("I9").Value = 2
$ws.Range("I10").Value = 2
$ws.Range("I11").Value = 2

$ws.Range("O6").Value = "kies…"
$ws.Range("O7").Value = "kies…"
$ws.Range("O8").Value = "kies…"
$ws.Range("O9").Value = "kies…"
$ws.Range("O10").Value = "kies…"
$ws.Range("O11").Value = "kies…"
$ws.Range("O18").Value = "kies…"
$ws.Range("O19").Value = "ja"
$ws.Range("O21").Value = "ja"
$ws.Range("O22").Value = "kies…"

# ---------------------------------------------------------------------
# "A 2020"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("A 2020")

$ws.Range("H9").Value = "Hoofdstuk 4 (Machtsfuncties) + Hoofdstuk 5 (Exponentiële functies) + Vaardigheden"

$ws.Range("O6").Value = "kies..."
$ws.Range("O7").Value = "kies..."
$ws.Range("O8").Value = "kies..."
$ws.Range("O9").Value = "kies..."
$ws.Range("O10").Value = "kies..."
$ws.Range("O11").Value = "kies..."
$ws.Range("O18").Value = "kies..."
$ws.Range("O19").Value = "ja"
$ws.Range("N20").Value = 1
$ws.Range("O21").Value = "ja"
$ws.Range("O22").Value = "kies..."

# ---------------------------------------------------------------------
# "A 2019"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("A 2019")

$ws.Range("O18").Value = "kies..."
$ws.Range("O19").Value = "ja"
$ws.Range("O21").Value = "ja"
$ws.Range("O22").Value = "kies..."

$ws.Range("O30").Value = "ja"
$ws.Range("N31").Value = 6
$ws.Range("O31").Value = "ja"
$ws.Range("O32").Value = "ja"

# ---------------------------------------------------------------------
# "A 2018"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("A 2018")

$ws.Range("O30").Value = "ja"
$ws.Range("O31").Value = "ja"
$ws.Range("O32").Value = "ja"

# Force a full recalculation so every dependent formula (AB/AF column
# flags, the F2 totals, the header summary text, …) picks up the new
# cached values.
$excel.CalculateFullRebuild()
